# Auto-generated script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.598.29"
$ws.Range("E2").Value = "  -4.94%  "
$ws.Range("D3").Value = "3.161.24"
$ws.Range("E3").Value = "  -5.58%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.32%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "3.139.57"
$ws.Range("E8").Value = "  -6.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.450"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.76%  "
$ws.Range("E11").Value = "  -8.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.77%  "
$ws.Range("D13").Value = "3.700.28"
$ws.Range("E13").Value = "  -5.52%  "
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.59%  "
$ws.Range("D16").Value = "3.168.78"
$ws.Range("E16").Value = "  -5.62%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "57.842.37"
$ws.Range("E17").Value = "  -4.54%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000154"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "350.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.08%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.512"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.71%  "
$ws.Range("D26").Value = "3.302.15"
$ws.Range("E26").Value = "  -5.47%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "0.0₃0960"
$ws.Range("E27").Value = "  -11.31%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.167"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.64%  "
$ws.Range("E35").Value = "  -6.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.32%  "
$ws.Range("E38").Value = "  -8.12%  "
$ws.Range("E39").Value = "  -8.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0700"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.58%  "
$ws.Range("D42").Value = "3.185.26"
$ws.Range("E42").Value = "  -5.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.87%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.695"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.29%  "
$ws.Range("E46").Value = "  -6.99%  "
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.57%  "
$ws.Range("D49").Value = "2.257.33"
$ws.Range("E49").Value = "  -8.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.38%  "
